# Implement csv module error handling
# Append a new data row (row 52) to each of the four log sheets, mirroring
# the structure of the existing rows (time stamp, hex byte strings, and
# decimal-decoded values).

$wb = $excel.ActiveWorkbook

$newRows = @{
    "ROW35-FE-LIFTER" = @{
        A = [double]"45750.35207362269"
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x72"
        E = "0xd"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 370
        I = 13
    }
    "ROW35-MID-LIFTER" = @{
        A = [double]"45750.20369950232"
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x72"
        E = "0xe"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 370
        I = 14
    }
    "ROW02-FE-LIFTER" = @{
        A = [double]"45750.34597035879"
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x72"
        E = "0x3"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 370
        I = 3
    }
    "ROW02-MID-LIFTER" = @{
        A = [double]"45750.40240554398"
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x72"
        E = "0x3"
        F = 400
        G = [double]"9.85046333984776e+23"
        H = 370
        I = 3
    }
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($newRows.ContainsKey($name)) {
        $data = $newRows[$name]
        $row = 52

        $ws.Cells.Item($row, 1).Value = $data.A
        $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
        $ws.Cells.Item($row, 2).Value = $data.B
        $ws.Cells.Item($row, 3).Value = $data.C
        $ws.Cells.Item($row, 4).Value = $data.D
        $ws.Cells.Item($row, 5).Value = $data.E
        $ws.Cells.Item($row, 6).Value = $data.F
        $ws.Cells.Item($row, 7).Value = $data.G
        $ws.Cells.Item($row, 8).Value = $data.H
        $ws.Cells.Item($row, 9).Value = $data.I
    }
}
